$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '''30.448.58'
$ws.Range("E2").Value = '  +1.09%  '
$ws.Range("D3").Value = '''1.854.48'
$ws.Range("E3").Value = '  +1.40%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''233.56'
$ws.Range("E5").Value = '  +1.07%  '
$ws.Range("D7").Value = '''0.4747'
$ws.Range("E7").Value = '  +1.96%  '
$ws.Range("D8").Value = '''0.2761'
$ws.Range("E8").Value = '  +2.91%  '
$ws.Range("D9").Value = '''0.06353'
$ws.Range("E9").Value = '  +1.51%  '
$ws.Range("D10").Value = '''17.98'
$ws.Range("E10").Value = '  +12.52%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = '''0.07477'
$ws.Range("E11").Value = '  +1.34%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '''1.820.64'
$ws.Range("E12").Value = '  -0.46%  '
$ws.Range("D13").Value = '''4.977'
$ws.Range("E13").Value = '  +1.66%  '
$ws.Range("E14").Value = '  +2.00%  '
$ws.Range("D15").Value = '''0.6253'
$ws.Range("E15").Value = '  +1.38%  '
$ws.Range("D16").Value = '''30.408.78'
$ws.Range("E16").Value = '  +1.19%  '
$ws.Range("D17").Value = '''246.07'
$ws.Range("E17").Value = '  +8.87%  '
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("D19").Value = '''12.71'
$ws.Range("E19").Value = '  +2.72%  '
$ws.Range("D20").Value = '''0.000007355'
$ws.Range("E20").Value = '  +1.18%  '
$ws.Range("D21").Value = '''0.9998'
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("D22").Value = '''4.939'
$ws.Range("E22").Value = '  +2.13%  '
$ws.Range("D23").Value = '''5.916'
$ws.Range("E23").Value = '  +1.02%  '
$ws.Range("D24").Value = '''164.34'
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").Value = '''9.076'
$ws.Range("E25").Value = '  -0.42%  '
$ws.Range("D26").Value = '''18.00'
$ws.Range("E26").Value = '  +2.30%  '
$ws.Range("D27").Value = '''1.884'
$ws.Range("E27").Value = '  +2.19%  '
$ws.Range("E28").Value = '  +1.52%  '
$ws.Range("D29").Value = '''1.346'
$ws.Range("E29").Value = '  -1.56%  '
$ws.Range("D30").Value = '''4.058'
$ws.Range("E30").Value = '  +0.39%  '
$ws.Range("E31").Value = '  +2.29%  '
$ws.Range("D32").Value = '''0.04847'
$ws.Range("E32").Value = '  +1.65%  '
$ws.Range("E33").Value = '  +0.92%  '
$ws.Range("D34").Value = '''0.6990'
$ws.Range("E34").Value = '  -0.51%  '
$ws.Range("D35").Value = '''2.699'
$ws.Range("E35").Value = '  +0.45%  '
$ws.Range("D36").Value = '''0.01900'
$ws.Range("E36").Value = '  +5.13%  '
$ws.Range("E37").Value = '  +3.01%  '
$ws.Range("D38").Value = '''0.8795'
$ws.Range("E38").Value = '  -1.21%  '
$ws.Range("D39").Value = '''1.994'
$ws.Range("E39").Value = '  +3.64%  '
$ws.Range("D40").Value = '''106.67'
$ws.Range("E40").Value = '  +3.46%  '
$ws.Range("D41").Value = '''0.9998'
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").Value = '''0.4077'
$ws.Range("E42").Value = '  +2.12%  '
$ws.Range("D43").Value = '''5.519'
$ws.Range("E43").Value = '  +1.07%  '
$ws.Range("D44").Value = '''7.194'
$ws.Range("E44").Value = '  +3.47%  '
$ws.Range("D45").Value = '''63.40'
$ws.Range("E45").Value = '  +6.46%  '
$ws.Range("E46").Value = '  +1.28%  '
$ws.Range("D47").Value = '''34.03'
$ws.Range("E47").Value = '  +4.37%  '
$ws.Range("D48").Value = '''8.552'
$ws.Range("E48").Value = '  +1.51%  '
$ws.Range("D49").Value = '''0.05500'
$ws.Range("E49").Value = '  -0.24%  '
$ws.Range("E50").Value = '  -0.83%  '
$ws.Range("D51").Value = '''0.3697'
$ws.Range("E51").Value = '  +2.40%  '
